$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sprint report date (B2): 2018-09-19 -> 2018-12-19
$ws.Range("B2").Value = 43453

# Clear the previous sprint's free-text answers that no longer apply to
# this sprint's report. Clearing removes both the cell value and its
# shared-string reference, leaving the cell blank but keeping its style.
$ws.Range("A10").Value = $null
$ws.Range("B10").Value = $null
$ws.Range("A11").Value = $null
$ws.Range("B11").Value = $null
$ws.Range("A15").Value = $null
$ws.Range("B15").Value = $null
$ws.Range("A16").Value = $null
$ws.Range("B16").Value = $null
$ws.Range("A17").Value = $null
$ws.Range("B17").Value = $null
$ws.Range("A25").Value = $null
$ws.Range("B25").Value = $null
$ws.Range("A28").Value = $null

# Match the author's final on-sheet selection/cursor position.
$ws.Range("G25").Select() | Out-Null
